# Add support for a ControlStrategy (curtailment) for Producers
# on the ConsumerProducer sheet: insert two new columns (Control_Strategy,
# CS_Curt_MaxPower) right before the existing MC_Value / MC_Profile_ID columns.

$wb = $excel.ActiveWorkbook

$wsCP = $wb.Worksheets.Item("ConsumerProducer")

# Insert two new blank columns at U (21) and V (22), pushing the existing
# MC_Value / MC_Profile_ID columns to W and X.
$wsCP.Columns.Item(21).Resize(1, 2).EntireColumn.Insert()

# Remove the leftover selection so the saved view just shows the default
# top-left selection (matches target: <sheetView workbookViewId="0"/>).
$wsCP.Cells.Item(1, 1).Select()

$wsCP.Cells.Item(1, 21).Value = "Control_Strategy"
$wsCP.Cells.Item(1, 22).Value = "CS_Curt_MaxPower"

# Give the two new columns the same "best fit" widths Excel would have
# computed for these header strings (16 and ~18.57 characters wide).
$wsCP.Columns.Item(21).ColumnWidth = 15.166666666666666
$wsCP.Columns.Item(22).ColumnWidth = 17.7

# Switch the active sheet from Conversion (tab index 4) to Areas (tab index 0),
# so "Conversion" no longer is the one marked as tab-selected.
$wsAreas = $wb.Worksheets.Item("Areas")
$wsAreas.Select()
$wsAreas.Cells.Item(1, 1).Select()
